# Auto-generated: update market-price snapshot cells per scheduled-runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 364.62
$ws.Range("J17").Value = 364.62
$ws.Range("L17").Value = 1093.86
$ws.Range("N17").Value = -1429.86
$ws.Range("H28").Value = 618.7778
$ws.Range("I28").Value = 618.7778
$ws.Range("K28").Value = 618.7778
$ws.Range("M28").Value = -133.7778
$ws.Range("H33").Value = 623.8095
$ws.Range("I33").Value = 662.5
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 662.5
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = -433.5
$ws.Range("N33").Value = -958
$ws.Range("H64").Value = 5009.1055
$ws.Range("I64").Value = 4724.875
$ws.Range("J64").Value = 5215.8184
$ws.Range("K64").Value = 4724.875
$ws.Range("L64").Value = 5215.8184
$ws.Range("M64").Value = -4476.875
$ws.Range("N64").Value = -5711.8184
$ws.Range("H67").Value = 5009.1055
$ws.Range("I67").Value = 4724.875
$ws.Range("J67").Value = 5215.8184
$ws.Range("K67").Value = 4724.875
$ws.Range("L67").Value = 5215.8184
$ws.Range("M67").Value = -3866.875
$ws.Range("N67").Value = -6931.8184
$ws.Range("H88").Value = 7633.3335
$ws.Range("I88").Value = 3866.6667
$ws.Range("K88").Value = 3866.6667
$ws.Range("M88").Value = -3460.6667
$ws.Range("H91").Value = 7633.3335
$ws.Range("I91").Value = 3866.6667
$ws.Range("K91").Value = 3866.6667
$ws.Range("M91").Value = -2462.6667
$ws.Range("H116").Value = 98853.37
$ws.Range("J116").Value = 3160
$ws.Range("L116").Value = 3160
$ws.Range("N116").Value = -10044
$ws.Range("H118").Value = 334696.66
$ws.Range("I118").Value = 500595
$ws.Range("J118").Value = 2900
$ws.Range("K118").Value = 1501785
$ws.Range("L118").Value = 8700
$ws.Range("M118").Value = -1500128
$ws.Range("N118").Value = -12014
$ws.Range("H132").Value = 3181
$ws.Range("I132").Value = 1478.8049
$ws.Range("J132").Value = 10160
$ws.Range("K132").Value = 4436.4147
$ws.Range("L132").Value = 30480
$ws.Range("M132").Value = -1906.4147
$ws.Range("N132").Value = -35540
$ws.Range("H137").Value = 5422.514
$ws.Range("I137").Value = 6288.7036
$ws.Range("J137").Value = 2499.125
$ws.Range("K137").Value = 18866.1108
$ws.Range("L137").Value = 7497.375
$ws.Range("M137").Value = -16316.1108
$ws.Range("N137").Value = -12597.375
$ws.Range("H141").Value = 3555.1428
$ws.Range("I141").Value = 2000
$ws.Range("J141").Value = 3814.3333
$ws.Range("K141").Value = 6000
$ws.Range("L141").Value = 11442.9999
$ws.Range("M141").Value = -820
$ws.Range("N141").Value = -21802.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1214
$ws.Range("I45").Value = 776.375
$ws.Range("J45").Value = 1714.1428
$ws.Range("K45").Value = 776.375
$ws.Range("L45").Value = 1714.1428
$ws.Range("M45").Value = -399.375
$ws.Range("N45").Value = -2468.1428
$ws.Range("H61").Value = 3991.3794
$ws.Range("I61").Value = 3858.6155
$ws.Range("J61").Value = 4099.25
$ws.Range("K61").Value = 3858.6155
$ws.Range("L61").Value = 4099.25
$ws.Range("M61").Value = -3646.6155
$ws.Range("N61").Value = -4523.25
$ws.Range("H74").Value = 1689.7222
$ws.Range("I74").Value = 1032.175
$ws.Range("J74").Value = 3568.4285
$ws.Range("K74").Value = 1032.175
$ws.Range("L74").Value = 3568.4285
$ws.Range("M74").Value = -158.175
$ws.Range("N74").Value = -5316.4285
$ws.Range("H77").Value = 1689.7222
$ws.Range("I77").Value = 1032.175
$ws.Range("J77").Value = 3568.4285
$ws.Range("K77").Value = 5160.875
$ws.Range("L77").Value = 17842.1425
$ws.Range("M77").Value = -792.875
$ws.Range("N77").Value = -26578.1425
$ws.Range("H136").Value = 3991.3794
$ws.Range("I136").Value = 3858.6155
$ws.Range("J136").Value = 4099.25
$ws.Range("K136").Value = 11575.8465
$ws.Range("L136").Value = 12297.75
$ws.Range("M136").Value = -9025.8465
$ws.Range("N136").Value = -17397.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1179.8182
$ws.Range("I94").Value = 1240
$ws.Range("J94").Value = 1129.6666
$ws.Range("K94").Value = 1240
$ws.Range("L94").Value = 1129.6666
$ws.Range("M94").Value = -789
$ws.Range("N94").Value = -2031.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1892.8254
$ws.Range("I31").Value = 1180.48
$ws.Range("J31").Value = 4632.615
$ws.Range("K31").Value = 1180.48
$ws.Range("L31").Value = 4632.615
$ws.Range("M31").Value = -885.48
$ws.Range("N31").Value = -5222.615
$ws.Range("H34").Value = 1892.8254
$ws.Range("I34").Value = 1180.48
$ws.Range("J34").Value = 4632.615
$ws.Range("K34").Value = 1180.48
$ws.Range("L34").Value = 4632.615
$ws.Range("M34").Value = -978.48
$ws.Range("N34").Value = -5036.615
$ws.Range("H58").Value = 1762.9166
$ws.Range("I58").Value = 1107.8572
$ws.Range("J58").Value = 2680
$ws.Range("K58").Value = 1107.8572
$ws.Range("L58").Value = 2680
$ws.Range("M58").Value = -904.8571999999999
$ws.Range("N58").Value = -3086
$ws.Range("H99").Value = 80400.84
$ws.Range("I99").Value = 252153
$ws.Range("J99").Value = 4066.5557
$ws.Range("K99").Value = 252153
$ws.Range("L99").Value = 4066.5557
$ws.Range("M99").Value = -250655
$ws.Range("N99").Value = -7062.5557
$ws.Range("H105").Value = 922.8525
$ws.Range("I105").Value = 837.5685999999999
$ws.Range("J105").Value = 1357.8
$ws.Range("K105").Value = 837.5685999999999
$ws.Range("L105").Value = 1357.8
$ws.Range("M105").Value = 909.4314000000001
$ws.Range("N105").Value = -4851.8
$ws.Range("H126").Value = 80400.84
$ws.Range("I126").Value = 252153
$ws.Range("J126").Value = 4066.5557
$ws.Range("K126").Value = 756459
$ws.Range("L126").Value = 12199.6671
$ws.Range("M126").Value = -753989
$ws.Range("N126").Value = -17139.6671
$ws.Range("H132").Value = 2037.1163
$ws.Range("I132").Value = 996.13794
$ws.Range("J132").Value = 4193.4287
$ws.Range("K132").Value = 2988.41382
$ws.Range("L132").Value = 12580.2861
$ws.Range("M132").Value = -458.4138199999998
$ws.Range("N132").Value = -17640.2861
$ws.Range("H134").Value = 1732.1154
$ws.Range("I134").Value = 912.1667
$ws.Range("K134").Value = 2736.5001
$ws.Range("M134").Value = -201.5001000000002
$ws.Range("H136").Value = 1762.9166
$ws.Range("I136").Value = 1107.8572
$ws.Range("J136").Value = 2680
$ws.Range("K136").Value = 3323.5716
$ws.Range("L136").Value = 8040
$ws.Range("M136").Value = -773.5715999999998
$ws.Range("N136").Value = -13140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 766.5454999999999
$ws.Range("I44").Value = 714.6667
$ws.Range("J44").Value = 1000
$ws.Range("K44").Value = 2144.0001
$ws.Range("L44").Value = 3000
$ws.Range("M44").Value = -1746.0001
$ws.Range("N44").Value = -3796
$ws.Range("H129").Value = 1385.5714
$ws.Range("I129").Value = 776.6667
$ws.Range("J129").Value = 1842.25
$ws.Range("K129").Value = 2330.0001
$ws.Range("L129").Value = 5526.75
$ws.Range("M129").Value = 2669.9999
$ws.Range("N129").Value = -15526.75
$ws.Range("H137").Value = 5741.8667
$ws.Range("I137").Value = 1716.0555
$ws.Range("J137").Value = 11780.583
$ws.Range("K137").Value = 5148.166499999999
$ws.Range("L137").Value = 35341.749
$ws.Range("M137").Value = -48.16649999999936
$ws.Range("N137").Value = -45541.749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1072.5
$ws.Range("I97").Value = 1072.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1072.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -576.5
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 5086.421
$ws.Range("I102").Value = 3625.2307
$ws.Range("J102").Value = 8252.333000000001
$ws.Range("K102").Value = 3625.2307
$ws.Range("L102").Value = 8252.333000000001
$ws.Range("M102").Value = -2003.2307
$ws.Range("N102").Value = -11496.333
$ws.Range("H132").Value = 4762.385
$ws.Range("I132").Value = 4569.684
$ws.Range("J132").Value = 5285.4287
$ws.Range("K132").Value = 13709.052
$ws.Range("L132").Value = 15856.2861
$ws.Range("M132").Value = -11179.052
$ws.Range("N132").Value = -20916.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 35000
$ws.Range("J115").Value = 35000
$ws.Range("L115").Value = 35000
$ws.Range("N115").Value = -37350
$ws.Range("H132").Value = 14714335
$ws.Range("I132").Value = 25002680
$ws.Range("J132").Value = 16699.857
$ws.Range("K132").Value = 75008040
$ws.Range("L132").Value = 50099.571
$ws.Range("M132").Value = -75005510
$ws.Range("N132").Value = -55159.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21741114
$ws.Range("I132").Value = 32259202
$ws.Range("J132").Value = 3729.5334
$ws.Range("K132").Value = 96777606
$ws.Range("M132").Value = -96775076
$ws.Range("N132").Value = -16248.6002
